# Apply the "tags" column (M) on the QAGlist_Teil1 sheet, and restore the
# active-sheet / selection state to match the post-edit workbook.
#
# Commit: "bugfixes: import auf 'kinderbetreuung' (auch für das logging),
# plots.R ebenfalls sourcen, examples für plot_summary_teil_1 vervollständigt,
# tags für alle fragen"
#
# The substantive spreadsheet change is the "tags für alle fragen" part:
# every question row (2-33) in QAGlist_Teil1 gets a tag value written into
# column M ("Gap1_type"), classifying the question as one (or more,
# comma-separated) of: Arbeit / Haushalt & Selbstsorge / Soziales Umfeld.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("QAGlist_Teil1")

# row -> tag text, in the same top-to-bottom order the original author
# filled them in (this also reproduces the shared-string insertion order:
# Arbeit, Haushalt & Selbstsorge, Soziales Umfeld, Arbeit,Haushalt & Selbstsorge)
$tags = @{
    2  = "Arbeit"
    3  = "Haushalt & Selbstsorge"
    4  = "Soziales Umfeld"
    5  = "Soziales Umfeld"
    6  = "Soziales Umfeld"
    7  = "Soziales Umfeld"
    8  = "Soziales Umfeld"
    9  = "Soziales Umfeld"
    10 = "Soziales Umfeld"
    11 = "Haushalt & Selbstsorge"
    12 = "Haushalt & Selbstsorge"
    13 = "Haushalt & Selbstsorge"
    14 = "Arbeit"
    15 = "Arbeit"
    16 = "Arbeit"
    17 = "Haushalt & Selbstsorge"
    18 = "Arbeit"
    19 = "Arbeit"
    20 = "Haushalt & Selbstsorge"
    21 = "Arbeit,Haushalt & Selbstsorge"
    22 = "Arbeit,Haushalt & Selbstsorge"
    23 = "Arbeit"
    24 = "Arbeit"
    25 = "Arbeit,Haushalt & Selbstsorge"
    26 = "Haushalt & Selbstsorge"
    27 = "Haushalt & Selbstsorge"
    28 = "Arbeit"
    29 = "Arbeit"
    30 = "Arbeit"
    31 = "Haushalt & Selbstsorge"
    32 = "Haushalt & Selbstsorge"
    33 = "Haushalt & Selbstsorge"
}

for ($row = 2; $row -le 33; $row++) {
    $ws1.Range("M$row").Value2 = $tags[$row]
}

# The workbook's active tab moves from Qlist_Teil2b (3rd sheet) back to
# QAGlist_Teil1 (1st sheet), and the last selection on QAGlist_Teil1 is
# M31 (the last cell touched while filling in the tags).
$ws1.Activate()
$ws1.Range("M31").Select()

$ws3 = $wb.Worksheets.Item("Qlist_Teil2b")
$ws3.Range("A7").Select()
